# Update nombre_aides (column C) and montant_total (column D) for the
# 2020-07-23 data refresh of the Fonds de solidarite volet 1 dataset.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @(
    @{Row=2; C=36099; D=52256297},
    @{Row=3; C=87758; D=128764218},
    @{Row=4; C=30077; D=44575318},
    @{Row=5; C=8341; D=12399676},
    @{Row=6; C=1805; D=2684106},
    @{Row=7; C=135; D=197593},
    @{Row=11; C=39441; D=53616874},
    @{Row=12; C=9280; D=13428808},
    @{Row=13; C=25163; D=36924088},
    @{Row=14; C=8050; D=11953896},
    @{Row=15; C=2071; D=3081650},
    @{Row=16; C=384; D=565123},
    @{Row=19; C=9767; D=12979510},
    @{Row=20; C=12944; D=18699771},
    @{Row=21; C=30780; D=45209323},
    @{Row=22; C=9974; D=14833195},
    @{Row=23; C=2542; D=3783263},
    @{Row=24; C=461; D=685845},
    @{Row=26; C=11243; D=15062605},
    @{Row=27; C=7355; D=10658114},
    @{Row=28; C=21801; D=32009957},
    @{Row=29; C=7569; D=11265756},
    @{Row=30; C=1885; D=2813044},
    @{Row=31; C=325; D=484915},
    @{Row=33; C=8000; D=10590622},
    @{Row=34; C=3002; D=4326152},
    @{Row=35; C=7382; D=10786529},
    @{Row=36; C=2974; D=4405288},
    @{Row=37; C=788; D=1174763},
    @{Row=38; C=143; D=212732},
    @{Row=40; C=2275; D=3069777},
    @{Row=41; C=16584; D=23994421},
    @{Row=42; C=49452; D=72538353},
    @{Row=43; C=18436; D=27388616},
    @{Row=44; C=5387; D=8024642},
    @{Row=45; C=1095; D=1633292},
    @{Row=46; C=57; D=83848},
    @{Row=49; C=16041; D=21399691},
    @{Row=50; C=1800; D=2612214},
    @{Row=51; C=6312; D=9288408},
    @{Row=52; C=2185; D=3263350},
    @{Row=53; C=713; D=1064805},
    @{Row=54; C=164; D=242833},
    @{Row=56; C=5879; D=8114927},
    @{Row=57; C=759; D=1113990},
    @{Row=58; C=1945; D=2885326},
    @{Row=59; C=795; D=1184953},
    @{Row=60; C=264; D=395758},
    @{Row=61; C=61; D=91500},
    @{Row=62; C=15; D=22500},
    @{Row=63; C=1122; D=1589709},
    @{Row=64; C=14786; D=21374285},
    @{Row=65; C=43420; D=63574400},
    @{Row=66; C=15266; D=22696732},
    @{Row=67; C=4415; D=6576560},
    @{Row=68; C=872; D=1297596},
    @{Row=71; C=14572; D=19261616},
    @{Row=72; C=48233; D=70234148},
    @{Row=73; C=138493; D=204137310},
    @{Row=74; C=60542; D=90245534},
    @{Row=75; C=19227; D=28733580},
    @{Row=76; C=4377; D=6540022},
    @{Row=77; C=234; D=346170},
    @{Row=78; C=19; D=27405},
    @{Row=83; C=47766; D=65251590},
    @{Row=84; C=4338; D=6290302},
    @{Row=85; C=11023; D=16198651},
    @{Row=86; C=3752; D=5592187},
    @{Row=87; C=1305; D=1949989},
    @{Row=88; C=276; D=411512},
    @{Row=89; C=22; D=32902},
    @{Row=91; C=5031; D=6781899},
    @{Row=92; C=1460; D=2111280},
    @{Row=93; C=4816; D=7094906},
    @{Row=94; C=1836; D=2735899},
    @{Row=95; C=653; D=978641},
    @{Row=96; C=165; D=246613},
    @{Row=99; C=3223; D=4276971},
    @{Row=100; C=555; D=827964},
    @{Row=101; C=306; D=457130},
    @{Row=102; C=110; D=165000},
    @{Row=103; C=40; D=60000},
    @{Row=104; C=20; D=30000},
    @{Row=105; C=10374; D=15069780},
    @{Row=106; C=28416; D=41764120},
    @{Row=107; C=9522; D=14161153},
    @{Row=108; C=2589; D=3860749},
    @{Row=109; C=456; D=681482},
    @{Row=110; C=45; D=67500},
    @{Row=112; C=9424; D=12475867},
    @{Row=113; C=29233; D=42196622},
    @{Row=114; C=64192; D=93995087},
    @{Row=115; C=20782; D=30894979},
    @{Row=116; C=5820; D=8671640},
    @{Row=117; C=1060; D=1584506},
    @{Row=118; C=66; D=96420},
    @{Row=121; C=24792; D=33175377},
    @{Row=122; C=34219; D=49441956},
    @{Row=123; C=73808; D=108021366},
    @{Row=124; C=22984; D=34123241},
    @{Row=125; C=6118; D=9095777},
    @{Row=126; C=1142; D=1698319},
    @{Row=130; C=30171; D=40163109},
    @{Row=131; C=12767; D=18487164},
    @{Row=132; C=31406; D=46150210},
    @{Row=133; C=11181; D=16613059},
    @{Row=134; C=2856; D=4258581},
    @{Row=135; C=458; D=680990},
    @{Row=136; C=31; D=46039},
    @{Row=138; C=10438; D=13953104},
    @{Row=139; C=33464; D=48367362},
    @{Row=140; C=78250; D=114700384},
    @{Row=141; C=23514; D=34956625},
    @{Row=142; C=6124; D=9140999},
    @{Row=143; C=1352; D=2012099},
    @{Row=144; C=71; D=106130},
    @{Row=145; C=13; D=19500},
    @{Row=146; C=27956; D=37840034}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    $ws.Cells.Item($u.Row, 4).Value = $u.D
}

$wb.Save()